# ParkingLot workbook: add Story 2 / Story 3 test-case sheets, fix a typo on Story 1.

$wb = $excel.ActiveWorkbook
$story1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Story 1: fix "cards" -> "cars" typo, and update the remembered selection.
# ---------------------------------------------------------------------------
$story1.Range("B10").Value = "Have 2 cars on the space"
$story1.Range("A6:F6").Select()

# ---------------------------------------------------------------------------
# Story 2: new sheet, formatted like Story 1, filled with the AC1/AC2 cases.
# ---------------------------------------------------------------------------
$story2 = $wb.Worksheets.Add($null, $story1)
$story2.Name = "Story 2"

$story1.Range("A2:G10").Copy() | Out-Null
$story2.Range("A2:G10").PasteSpecial(-4122) | Out-Null

$story2.Range("A2").Value = "STORY"
$story2.Range("B2").Value = "2"
$story2.Range("B3").Value = "GIVEN"

$story2.Range("A4").Value = "Test Case"
$story2.Range("B4").Value = "Parking Lot Status"
$story2.Range("C4").Value = "Car"
$story2.Range("D4").Value = "Ticket"
$story2.Range("E4").Value = "When"
$story2.Range("F4").Value = "Then"
$story2.Range("G4").Value = "Remarks"

$story2.Range("A5").Value = "Fetch Car (AC1)"
$story2.Range("B5").Value = "have cars"
$story2.Range("C5").Value = "does not matter"
$story2.Range("D5").Value = "Wrong Ticket"
$story2.Range("E5").Value = "Fetch Car"
$story2.Range("F5").Value = "Will not return the car and display message of ""Unrecognized parking ticket"""

$story2.Range("A6").Value = "Fetch Car (AC1)"
$story2.Range("B6").Value = "have cars"
$story2.Range("C6").Value = "already fetched"
$story2.Range("D6").Value = "used ticket"
$story2.Range("E6").Value = "Fetch Car"
$story2.Range("F6").Value = "Will not return the car and display message of ""Unrecognized parking ticket"""

$story2.Range("A7").Value = "Park Car (AC2)"
$story2.Range("B7").Value = "NO open space"
$story2.Range("C7").Value = "existing"
$story2.Range("D7").Value = "-"
$story2.Range("E7").Value = "Park Car"
$story2.Range("F7").Value = "Will not return a ticket and display message of ""No available position"""

$story2.Range("B3:D3").Merge()

$story2.Rows.Item(5).RowHeight = 30
$story2.Rows.Item(6).RowHeight = 30
$story2.Rows.Item(7).RowHeight = 30

$story2.Columns.Item(1).ColumnWidth = 15.92
$story2.Columns.Item(2).ColumnWidth = 17.08
$story2.Columns.Item(3).ColumnWidth = 19.58
$story2.Columns.Item(4).ColumnWidth = 14.08
$story2.Columns.Item(5).ColumnWidth = 11.08
$story2.Columns.Item(6).ColumnWidth = 39.42
$story2.Columns.Item(7).ColumnWidth = 25.75
$story2.Columns.Item(8).ColumnWidth = 13.42

# ---------------------------------------------------------------------------
# Story 3: new sheet, formatted like Story 1, only the header rows filled in.
# ---------------------------------------------------------------------------
$story3 = $wb.Worksheets.Add($null, $story2)
$story3.Name = "Story 3"

$story1.Range("A2:G10").Copy() | Out-Null
$story3.Range("A2:G10").PasteSpecial(-4122) | Out-Null

$story3.Range("A2").Value = "STORY"
$story3.Range("B2").Value = "3"
$story3.Range("B3").Value = "GIVEN"

$story3.Range("A4").Value = "Test Case"
$story3.Range("B4").Value = "Parking Lot Status"
$story3.Range("C4").Value = "Car"
$story3.Range("D4").Value = "Ticket"
$story3.Range("E4").Value = "When"
$story3.Range("F4").Value = "Then"
$story3.Range("G4").Value = "Remarks"

$story3.Range("B3:D3").Merge()

$story3.Range("A5:G10").ClearContents()
$story3.Range("A8:G10").EntireRow.Delete()

$story3.Columns.Item(1).ColumnWidth = 15.58
$story3.Columns.Item(2).ColumnWidth = 16.58
$story3.Columns.Item(3).ColumnWidth = 20.25
$story3.Columns.Item(4).ColumnWidth = 19.58
$story3.Columns.Item(5).ColumnWidth = 15.25
$story3.Columns.Item(6).ColumnWidth = 41.08

$story3.Range("A5").Select()

# ---------------------------------------------------------------------------
# Selections / active sheet: Story 2 ends up the active tab, like in the diff.
# ---------------------------------------------------------------------------
$story2.Range("D9").Select()
$story2.Activate()
